$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 97.73
$ws.Range("K2").Value = 11.36
$ws.Range("N2").Value = 9.09
$ws.Range("G3").Value = 98.7
$ws.Range("K3").Value = 1.3
$ws.Range("N3").Value = 3.9
$ws.Range("G4").Value = 96.3
$ws.Range("K4").Value = 7.41
$ws.Range("N4").Value = 11.11
$ws.Range("K5").Value = 11.54
$ws.Range("N5").Value = 11.54
$ws.Range("K6").Value = 9.94
$ws.Range("N6").Value = 7.18
$ws.Range("K7").Value = 13.04
$ws.Range("N7").Value = 9.57
$ws.Range("G8").Value = 98.88
$ws.Range("K8").Value = 21.01
$ws.Range("N8").Value = 16.25
$ws.Range("K9").Value = 13.16
$ws.Range("N9").Value = 10.53
$ws.Range("G10").Value = 99.06999999999999
$ws.Range("K10").Value = 9.720000000000001
$ws.Range("N10").Value = 8.800000000000001
$ws.Range("K11").Value = 4.55
$ws.Range("N11").Value = 4.55
$ws.Range("G12").Value = 99.34
$ws.Range("K12").Value = 8.94
$ws.Range("N12").Value = 5.3
$ws.Range("K13").Value = 7.55
$ws.Range("N13").Value = 5.66
$ws.Range("K14").Value = 11.36
$ws.Range("N14").Value = 6.82
$ws.Range("K15").Value = 6.21
$ws.Range("N15").Value = 8.07
$ws.Range("K16").Value = 10.81
$ws.Range("N16").Value = 10.81
$ws.Range("K17").Value = 10.53
$ws.Range("N17").Value = 10
$ws.Range("K18").Value = 3.31
$ws.Range("N18").Value = 3.31
$ws.Range("G19").Value = 99.16
$ws.Range("K19").Value = 6.75
$ws.Range("N19").Value = 5.91
$ws.Range("G20").Value = 99.73999999999999
$ws.Range("K20").Value = 6.23
$ws.Range("N20").Value = 4.68
$ws.Range("K21").Value = 14.06
$ws.Range("N21").Value = 10.94
$ws.Range("G22").Value = 99.29000000000001
$ws.Range("K22").Value = 2.84
$ws.Range("N22").Value = 2.13
$ws.Range("K23").Value = 10
$ws.Range("N23").Value = 10
$ws.Range("G24").Value = 97.73
$ws.Range("K24").Value = 9.09
$ws.Range("N24").Value = 6.82
$ws.Range("G25").Value = 99.43000000000001
$ws.Range("K25").Value = 6.25
$ws.Range("N25").Value = 5.68
$ws.Range("K26").Value = 13.28
$ws.Range("N26").Value = 12.5
$ws.Range("G27").Value = 99.64
$ws.Range("K27").Value = 4.69
$ws.Range("N27").Value = 3.97
$ws.Range("K28").Value = 17.5
$ws.Range("N28").Value = 17.5
